$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9803115129470825
$ws.Range("B1").Value = 1.120945930480957
$ws.Range("C1").Value = 1.40932834148407
$ws.Range("D1").Value = 2.887648582458496
$ws.Range("E1").Value = 4.316815853118896
